# Weekly update for Hortaliza, Terminal Hortofrutícola Agro Chillán - Berenjena
# A new observation row is inserted at row 19 (pushing the previous rows 19-22
# down to 20-23), and the new row is populated with the latest market data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 19, shifting rows 19-22 -> 20-23
$ws.Rows.Item(19).Insert()

# Copy the date cell's number format from the row below so the new date cell
# keeps the same date/time style used throughout column D.
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat

# Populate the new row 19 with the new record's values
$ws.Cells.Item(19, 1).Value = 7
$ws.Cells.Item(19, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19, 3).Value = "Ñuble"
$ws.Cells.Item(19, 4).Value = 44617
$ws.Cells.Item(19, 5).Value = 16
$ws.Cells.Item(19, 6).Value = 100112001
$ws.Cells.Item(19, 7).Value = "Berenjena"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 11).Value = 10000
$ws.Cells.Item(19, 12).Value = 11000
$ws.Cells.Item(19, 13).Value = 10500
$ws.Cells.Item(19, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(19, 15).Value = "Región Metropolitana"
$ws.Cells.Item(19, 16).Value = 175
$ws.Cells.Item(19, 17).Value = 60
$ws.Cells.Item(19, 18).Value = "Hortaliza"
